$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.218.27'
$ws.Range("E2").Value = '  -0.06%  '

# Row 3
$ws.Range("D3").Value = '1.854.77'
$ws.Range("E3").Value = '  -0.43%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.36'
$ws.Range("E5").Value = '  -0.34%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6969'
$ws.Range("E6").Value = '  -1.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07725'
$ws.Range("E8").Value = '  -1.12%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3070'
$ws.Range("E9").Value = '  -1.38%  '

# Row 10
$ws.Range("E10").Value = '  -2.19%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07808'

# Row 12
$ws.Range("D12").Value = '1.866.50'
$ws.Range("E12").Value = '  +0.29%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.096'
$ws.Range("E13").Value = '  -1.50%  '

# Row 14
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.12'
$ws.Range("E14").Value = '  -1.52%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6861'
$ws.Range("E15").Value = '  -1.36%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.518'
$ws.Range("E16").Value = '  +2.63%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008414'
$ws.Range("E17").Value = '  +1.36%  '

# Row 18
$ws.Range("D18").Value = '29.224.37'
$ws.Range("E18").Value = '  +0.11%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '249.02'
$ws.Range("E19").Value = '  -1.35%  '

# Row 20
$ws.Range("D20").Value = '2.111.90'
$ws.Range("E20").Value = '  -0.13%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.79'
$ws.Range("E21").Value = '  -2.25%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.10%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.505'
$ws.Range("E23").Value = '  +0.14%  '

# Row 24
$ws.Range("E24").Value = '  -0.13%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1500'
$ws.Range("E25").Value = '  -3.76%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.69'
$ws.Range("E26").Value = '  +0.81%  '

# Row 27
$ws.Range("E27").Value = '  -1.62%  '

# Row 28
$ws.Range("E28").Value = '  -1.90%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.562'
$ws.Range("E29").Value = '  +4.28%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.236'
$ws.Range("E30").Value = '  -1.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.179'
$ws.Range("E31").Value = '  -2.08%  '

# Row 32
$ws.Range("E32").Value = '  -1.30%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05206'
$ws.Range("E33").Value = '  -1.10%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7593'
$ws.Range("E34").Value = '  +2.06%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.841'
$ws.Range("E35").Value = '  -2.43%  '

# Row 36
$ws.Range("E36").Value = '  +0.73%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01859'
$ws.Range("E38").Value = '  +0.03%  '

# Row 39
$ws.Range("D39").Value = '1.213.63'
$ws.Range("E39").Value = '  -3.04%  '

# Row 40
$ws.Range("E40").Value = '  -0.75%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8968'
$ws.Range("E41").Value = '  -0.34%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.78'
$ws.Range("E42").Value = '  -1.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9993'
$ws.Range("E43").Value = '  -0.10%  '

# Row 44
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.521'
$ws.Range("E44").Value = '  -12.33%  '

# Row 45
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '2.012.09'
$ws.Range("E45").Value = '  -1.50%  '

# Row 46
$ws.Range("E46").Value = '  -3.40%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.20'
$ws.Range("E47").Value = '  -9.08%  '

# Row 48
$ws.Range("E48").Value = '  -0.38%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.530'
$ws.Range("E49").Value = '  +1.49%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.750'
$ws.Range("E50").Value = '  -1.62%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.003'
$ws.Range("E51").Value = '  +0.17%  '
